# Auto-generated script to apply scheduled-runner profit recalculation updates
# to the Ultros_Profits workbook's per-job Leve tables (columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2006.8462
$ws.Range("I70").Value = 1495.1666
$ws.Range("J70").Value = 2445.4285
$ws.Range("K70").Value = 4485.4998
$ws.Range("L70").Value = 7336.2855
$ws.Range("M70").Value = -4215.4998
$ws.Range("N70").Value = -7876.2855
$ws.Range("H73").Value = 2006.8462
$ws.Range("I73").Value = 1495.1666
$ws.Range("J73").Value = 2445.4285
$ws.Range("K73").Value = 4485.4998
$ws.Range("L73").Value = 7336.2855
$ws.Range("M73").Value = -3549.4998
$ws.Range("N73").Value = -9208.2855
$ws.Range("H86").Value = 2824.4
$ws.Range("I86").Value = 1972.5
$ws.Range("K86").Value = 1972.5
$ws.Range("M86").Value = -849.5
$ws.Range("H89").Value = 2824.4
$ws.Range("I89").Value = 1972.5
$ws.Range("K89").Value = 9862.5
$ws.Range("M89").Value = -4246.5
$ws.Range("H125").Value = 2684.5
$ws.Range("H127").Value = 8315.833000000001
$ws.Range("I127").Value = 1780
$ws.Range("J127").Value = 40995
$ws.Range("K127").Value = 5340
$ws.Range("L127").Value = 122985
$ws.Range("M127").Value = -380
$ws.Range("N127").Value = -132905
$ws.Range("H137").Value = 3376.842
$ws.Range("I137").Value = 3215.8823
$ws.Range("J137").Value = 4745
$ws.Range("K137").Value = 9647.6469
$ws.Range("L137").Value = 14235
$ws.Range("M137").Value = -7097.6469
$ws.Range("N137").Value = -19335
$ws.Range("H141").Value = 4527.593
$ws.Range("I141").Value = 2176.875
$ws.Range("K141").Value = 6530.625
$ws.Range("M141").Value = -1350.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15388742
$ws.Range("I32").Value = 15877213
$ws.Range("K32").Value = 15877213
$ws.Range("M32").Value = -15876926
$ws.Range("H61").Value = 1265.4166
$ws.Range("I61").Value = 1103.0435
$ws.Range("K61").Value = 1103.0435
$ws.Range("M61").Value = -891.0435
$ws.Range("H74").Value = 3068.077
$ws.Range("H77").Value = 3068.077
$ws.Range("H102").Value = 9901.637000000001
$ws.Range("I102").Value = 9790.700000000001
$ws.Range("K102").Value = 9790.700000000001
$ws.Range("M102").Value = -8168.700000000001
$ws.Range("H132").Value = 1556.3112
$ws.Range("I132").Value = 1642.439
$ws.Range("J132").Value = 673.5
$ws.Range("K132").Value = 4927.317
$ws.Range("L132").Value = 2020.5
$ws.Range("M132").Value = -2397.317
$ws.Range("N132").Value = -7080.5
$ws.Range("H136").Value = 1265.4166
$ws.Range("I136").Value = 1103.0435
$ws.Range("K136").Value = 3309.1305
$ws.Range("M136").Value = -759.1305000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2413.6956
$ws.Range("I86").Value = 2110.647
$ws.Range("J86").Value = 3272.3333
$ws.Range("K86").Value = 2110.647
$ws.Range("L86").Value = 3272.3333
$ws.Range("M86").Value = -987.6469999999999
$ws.Range("N86").Value = -5518.3333
$ws.Range("H89").Value = 2413.6956
$ws.Range("I89").Value = 2110.647
$ws.Range("J89").Value = 3272.3333
$ws.Range("K89").Value = 10553.235
$ws.Range("L89").Value = 16361.6665
$ws.Range("M89").Value = -4937.235000000001
$ws.Range("N89").Value = -27593.6665
$ws.Range("H97").Value = 8745.444
$ws.Range("I97").Value = 8745.444
$ws.Range("K97").Value = 8745.444
$ws.Range("M97").Value = -7754.444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 935.5
$ws.Range("I22").Value = 961.6667
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 961.6667
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -611.6667
$ws.Range("N22").Value = -1400
$ws.Range("H31").Value = 1579.4736
$ws.Range("I31").Value = 1568.1082
$ws.Range("K31").Value = 1568.1082
$ws.Range("M31").Value = -1273.1082
$ws.Range("H34").Value = 1579.4736
$ws.Range("I34").Value = 1568.1082
$ws.Range("K34").Value = 1568.1082
$ws.Range("M34").Value = -1366.1082
$ws.Range("H58").Value = 2342
$ws.Range("I58").Value = 1781.125
$ws.Range("K58").Value = 1781.125
$ws.Range("M58").Value = -1578.125
$ws.Range("H62").Value = 4000.625
$ws.Range("I62").Value = 4117.6665
$ws.Range("K62").Value = 4117.6665
$ws.Range("M62").Value = -3493.6665
$ws.Range("H65").Value = 4000.625
$ws.Range("I65").Value = 4117.6665
$ws.Range("K65").Value = 20588.3325
$ws.Range("M65").Value = -17468.3325
$ws.Range("H86").Value = 17795
$ws.Range("J86").Value = 9630.137000000001
$ws.Range("L86").Value = 9630.137000000001
$ws.Range("N86").Value = -11876.137
$ws.Range("H89").Value = 17795
$ws.Range("J89").Value = 9630.137000000001
$ws.Range("L89").Value = 48150.685
$ws.Range("N89").Value = -59382.685
$ws.Range("H132").Value = 3142.818
$ws.Range("I132").Value = 2957.1
$ws.Range("K132").Value = 8871.299999999999
$ws.Range("M132").Value = -6341.299999999999
$ws.Range("H136").Value = 2342
$ws.Range("I136").Value = 1781.125
$ws.Range("K136").Value = 5343.375
$ws.Range("M136").Value = -2793.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 613.6667
$ws.Range("I22").Value = 636.4
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1909.2
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -1740.2
$ws.Range("N22").Value = -1838
$ws.Range("H27").Value = 613.6667
$ws.Range("I27").Value = 636.4
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 1909.2
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -1807.2
$ws.Range("N27").Value = -1704
$ws.Range("H131").Value = 2858.1538
$ws.Range("J131").Value = 4351
$ws.Range("L131").Value = 13053
$ws.Range("N131").Value = -23133

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3275.724
$ws.Range("I102").Value = 2239.5557
$ws.Range("K102").Value = 2239.5557
$ws.Range("M102").Value = -617.5556999999999
$ws.Range("H126").Value = 7025.6
$ws.Range("I126").Value = 6966.6665
$ws.Range("J126").Value = 7114
$ws.Range("K126").Value = 20899.9995
$ws.Range("L126").Value = 21342
$ws.Range("M126").Value = -18429.9995
$ws.Range("N126").Value = -26282
$ws.Range("H132").Value = 2582.8235
$ws.Range("I132").Value = 2582.8235
$ws.Range("K132").Value = 7748.470499999999
$ws.Range("M132").Value = -5218.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 985.7037
$ws.Range("J22").Value = 1049.5652
$ws.Range("L22").Value = 1049.5652
$ws.Range("N22").Value = -1639.5652
$ws.Range("H27").Value = 985.7037
$ws.Range("J27").Value = 1049.5652
$ws.Range("L27").Value = 1049.5652
$ws.Range("N27").Value = -1263.5652
$ws.Range("H40").Value = 19727.666
$ws.Range("I40").Value = 22423.2
$ws.Range("K40").Value = 22423.2
$ws.Range("M40").Value = -22287.2
$ws.Range("H93").Value = 7455.3687
$ws.Range("J93").Value = 7652
$ws.Range("L93").Value = 7652
$ws.Range("N93").Value = -10148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14479.667
$ws.Range("I45").Value = 7967.3335
$ws.Range("K45").Value = 7967.3335
$ws.Range("M45").Value = -7476.3335
$ws.Range("H54").Value = 27499.75
$ws.Range("J54").Value = 27499.75
$ws.Range("L54").Value = 27499.75
$ws.Range("N54").Value = -28539.75
$ws.Range("H100").Value = 1750.3334
$ws.Range("I100").Value = 2310
$ws.Range("J100").Value = 631
$ws.Range("K100").Value = 4620
$ws.Range("L100").Value = 1262
$ws.Range("M100").Value = -4079
$ws.Range("N100").Value = -2344
$ws.Range("H132").Value = 2713.2856
$ws.Range("I132").Value = 2799.2
$ws.Range("K132").Value = 8397.599999999999
$ws.Range("M132").Value = -5867.599999999999
$ws.Range("H136").Value = 1804.1578
$ws.Range("I136").Value = 1372.3334
$ws.Range("J136").Value = 3423.5
$ws.Range("K136").Value = 4117.0002
$ws.Range("L136").Value = 10270.5
$ws.Range("M136").Value = -1567.0002
$ws.Range("N136").Value = -15370.5
